$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, shifting existing rows 100:128 down to 101:129
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with data (matches the template of surrounding rows)
$ws.Range("A100").Value = 1
$ws.Range("B100").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C100").Value = "Arica y Parinacota"
$ws.Range("D100").Value = 44711
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E100").Value = 15
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100108
$ws.Range("H100").Value = "Tropicales y subtropicales"
$ws.Range("I100").Value = 100108003
$ws.Range("J100").Value = "Maracuyá"
$ws.Range("K100").Value = "Sin especificar"
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 130
$ws.Range("N100").Value = 19000
$ws.Range("O100").Value = 20000
$ws.Range("P100").Value = 19500
$ws.Range("Q100").Value = "$/caja 20 kilos"
$ws.Range("R100").Value = "Región de Arica y Parinacota"
$ws.Range("S100").Value = 975
$ws.Range("T100").Value = 20
